$wb = $excel.ActiveWorkbook

# --- Update header on "Weekly Quantity" sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Update header on "Monthly Trend" sheet ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after "Monthly Trend" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/bordered header style used on the other sheets
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122) # xlPasteFormats

$data = @(
    @(44934.99999999999, 155, -48.87760969719855, 374.4188898096907),
    @(44941.99999999999, 157, -48.43149874507972, 380.9876198576562),
    @(44955.99999999999, 162, -59.87684536294487, 376.3364992586787),
    @(44969.99999999999, 166, -41.1779358107945, 385.5028246083195),
    @(44976.99999999999, 169, -49.87048728725262, 387.9274138215677),
    @(45137.99999999999, 221, 9.999341567356501, 446.515321702279),
    @(45144.99999999999, 223, 3.368584793364743, 431.2880681493549),
    @(45151.99999999999, 225, 9.495308301199803, 443.1435091689489),
    @(45158.99999999999, 227, 8.692578556583097, 444.11297509539),
    @(45165.99999999999, 230, 9.918021331882469, 449.3418928664794),
    @(45172.99999999999, 232, 17.0469357560508, 435.3401530201569),
    @(45179.99999999999, 234, 35.77039359254241, 457.7769181468025),
    @(45186.99999999999, 236, 24.70687686060238, 449.0834045304273),
    @(45193.99999999999, 239, 23.29969589447631, 436.8275592396438)
)

$rowIndex = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($rowIndex, 1).Value = $row[0]
    $wsForecast.Cells.Item($rowIndex, 2).Value = $row[1]
    $wsForecast.Cells.Item($rowIndex, 3).Value = $row[2]
    $wsForecast.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex++
}

# Match date-style formatting used for the "ds" date column on the other sheets
$wsForecast.Range("A2:A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
